$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated record_atd (C) and record_id (D) values for corrected relevance markers
$ws.Range("C3:D3").Value = 30
$ws.Range("C5:D5").Value = 70
$ws.Range("C7:D7").Value = 260
$ws.Range("C9:D9").Value = 42
$ws.Range("C11:D11").Value = 34
$ws.Range("C13:D13").Value = 18
$ws.Range("C15:D15").Value = 178
$ws.Range("C18:D18").Value = 16
$ws.Range("C19:D19").Value = 33
$ws.Range("C21:D21").Value = 256
$ws.Range("C23:D23").Value = 258
$ws.Range("C25:D25").Value = 264
$ws.Range("C27:D27").Value = 98
$ws.Range("C29:D29").Value = 31
$ws.Range("C31:D31").Value = 203
$ws.Range("C33:D33").Value = 19
$ws.Range("C35:D35").Value = 111
$ws.Range("C37:D37").Value = 144
$ws.Range("C38:D38").Value = 1294
$ws.Range("C40:D40").Value = 59
$ws.Range("C42:D42").Value = 332
$ws.Range("C44:D44").Value = 26
$ws.Range("C46:D46").Value = 204
$ws.Range("C48:D48").Value = 108
$ws.Range("C50:D50").Value = 93
$ws.Range("C52:D52").Value = 40
$ws.Range("C54:D54").Value = 172
$ws.Range("C56:D56").Value = 216
$ws.Range("C58:D58").Value = 119
$ws.Range("C60:D60").Value = 232
$ws.Range("C64:D64").Value = 20
$ws.Range("C66:D66").Value = 206
$ws.Range("C68:D68").Value = 270
$ws.Range("C70:D70").Value = 217
$ws.Range("C72:D72").Value = 212
$ws.Range("C74:D74").Value = 470
$ws.Range("C76:D76").Value = 218
$ws.Range("C78:D78").Value = 101
$ws.Range("C80:D80").Value = 21
$ws.Range("C82:D82").Value = 25
$ws.Range("C84:D84").Value = 223
$ws.Range("C86:D86").Value = 179
$ws.Range("C88:D88").Value = 104
$ws.Range("C90:D90").Value = 131
$ws.Range("C92:D92").Value = 495
$ws.Range("C94:D94").Value = 71
$ws.Range("C96:D96").Value = 306

# Recompute the average_simulation_TD summary cell (C97) to reflect the corrected values
$ws.Range("C97").Value = 172.2083333333333
